$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 data (supplier D12/Dhani/JL.merbabu -> GAC/Gacoan/Tlogomas)
$ws.Range("A2").Value = "GAC"
$ws.Range("B2").Value = "Gacoan"
$ws.Range("C2").Value = "Tlogomas"

# Remove row 3 (D20/Danan/JL.biohazard) entirely
$ws.Rows(3).Delete()

# Move the selection to F5 as in the final workbook
$ws.Range("F5").Select()
